$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "Test Ringover (NO TOCAR)" order row (row 4) as a second
# test row for order 69134d11b9c1d30b15fabdc3: copy the row, then insert
# the copied cells above themselves. This shifts the original row down to
# row 5 (keeping it byte-for-byte, including its unused placeholder
# cells) while row 4 receives a fresh paste of the same values.
$ws.Rows.Item(4).Copy()
$ws.Rows.Item(4).Insert()
